$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 19: "ejs" label in column A
$ws.Range("A19").Value = "ejs"

# New row 20: hyperlink text in column D, linking to the ejs includes docs
$ws.Hyperlinks.Add($ws.Range("D20"), "https://github.com/visionmedia/ejs", "includes", [Type]::Missing, "https://github.com/visionmedia/ejs - includes")
$ws.Range("D20").Value = "https://github.com/visionmedia/ejs#includes"
$ws.Range("D20").Style = "Hyperlink"

# Update the saved selection to the newly added cell
$ws.Range("B20").Select()
